$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.833.60"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'1.766.76"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'338.70"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("D8").Value = "'0.3362"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("D9").Value = "'45.75"
$ws.Range("E9").Value = "  -5.30%  "
$ws.Range("D10").Value = "'1.135"
$ws.Range("E10").Value = "  -5.54%  "
$ws.Range("D11").Value = "'0.07270"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").Value = "'22.98"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "'6.257"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'1.766.03"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "'0.06607"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "'81.24"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "'17.17"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").Value = "'6.339"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").Value = "'27.855.29"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'11.81"
$ws.Range("E24").Value = "  -8.29%  "
$ws.Range("D25").Value = "'2.375"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").Value = "'1.494"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'153.42"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'20.06"
$ws.Range("E28").Value = "  -5.96%  "
$ws.Range("D29").Value = "'2.360"
$ws.Range("E29").Value = "  -7.71%  "
$ws.Range("D30").Value = "'1.967.80"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "'132.48"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("D32").Value = "'4.036"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.924"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").Value = "'0.08755"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "'12.44"
$ws.Range("E35").Value = "  -6.52%  "
$ws.Range("D36").Value = "'0.02362"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'0.6713"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "'0.06273"
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").Value = "'5.212"
$ws.Range("E39").Value = "  -6.20%  "
$ws.Range("D40").Value = "'0.2122"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").Value = "'1.229"
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").Value = "'1.472"
$ws.Range("E42").Value = "  -8.75%  "
$ws.Range("D43").Value = "'8.084"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D45").Value = "'13.89"
$ws.Range("E45").Value = "  -6.37%  "
$ws.Range("D46").Value = "'0.6132"
$ws.Range("E46").Value = "  -6.30%  "
$ws.Range("D47").Value = "'3.841"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'131.60"
$ws.Range("D49").Value = "'2.029"
$ws.Range("E49").Value = "  -6.27%  "
$ws.Range("D50").Value = "'0.07282"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "'1.186"
$ws.Range("E51").Value = "  +1.66%  "
